# Logged Week 17 data and fixed Simulate_Season.py tiebreaking method
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Rushing sheet: insert a new row for backup QB C.Wentz (row 3),
# renumber the rank column (A) to be 0-indexed, and update the
# Week-17 rushing totals for the players whose stats changed.
# ---------------------------------------------------------------
$rush = $wb.Worksheets.Item("Rushing")

# Make room for the new C.Wentz row right after S.Ehlinger by shifting
# the existing player rows (3-10) down to (4-11), bottom-up so nothing
# gets clobbered. Using a full-row copy/paste (rather than Rows.Insert)
# keeps the original style indices intact instead of synthesising new
# blended ones.
for ($r = 10; $r -ge 3; $r--) {
    $dst = $r + 1
    $rush.Range("A" + $r + ":F" + $r).Copy()
    $rush.Range("A" + $dst + ":F" + $dst).PasteSpecial(-4104)
}
$rush.Range("A3:F3").ClearContents()

# Make sure the rank column keeps its original (bold/bordered) style all
# the way down, including the brand-new row 11.
$rush.Range("A2").Copy()
$rush.Range("A3:A11").PasteSpecial(-4122)

# Re-number the rank column A for every data row (now 0-based, 10 rows).
$rush.Range("A2").Value = 0
$rush.Range("A3").Value = 1
$rush.Range("A4").Value = 2
$rush.Range("A5").Value = 3
$rush.Range("A6").Value = 4
$rush.Range("A7").Value = 5
$rush.Range("A8").Value = 6
$rush.Range("A9").Value = 7
$rush.Range("A10").Value = 8
$rush.Range("A11").Value = 9

# New row: C.Wentz
$rush.Range("B3").Value = "C.Wentz"
$rush.Range("C3").Value = 0
$rush.Range("D3").Value = 1
$rush.Range("E3").Value = 2
$rush.Range("F3").Value = 0

# Updated Week 17 totals (rows shifted down by one after the insert).
# J.Taylor
$rush.Range("C4").Value = 172
$rush.Range("D4").Value = 112
$rush.Range("E4").Value = 32
$rush.Range("F4").Value = 80

# N.Hines
$rush.Range("C5").Value = 27
$rush.Range("D5").Value = 22
$rush.Range("E5").Value = 8
$rush.Range("F5").Value = 5

# De.Jackson
$rush.Range("C7").Value = 3
$rush.Range("D7").Value = 1
$rush.Range("E7").Value = 3
$rush.Range("F7").Value = 3

# ---------------------------------------------------------------
# Receiving sheet: update the Week-17 receiving totals. No rows
# are added/removed here, only values change.
# ---------------------------------------------------------------
$rec = $wb.Worksheets.Item("Receiving")

# J.Taylor (row 2)
$rec.Range("C2").Value = 45
$rec.Range("D2").Value = 33

# N.Hines (row 3)
$rec.Range("C3").Value = 48
$rec.Range("D3").Value = 35

# M.Pittman (row 5)
$rec.Range("C5").Value = 99
$rec.Range("D5").Value = 72
$rec.Range("G5").Value = 17
$rec.Range("H5").Value = 8

# Z.Pascal (row 6)
$rec.Range("C6").Value = 53
$rec.Range("D6").Value = 30
$rec.Range("E6").Value = 12

# A.Dulin (row 9)
$rec.Range("C9").Value = 16
$rec.Range("G9").Value = 2

# T.Hilton (row 10)
$rec.Range("C10").Value = 37
$rec.Range("E10").Value = 9
$rec.Range("F10").Value = 5

# D.Patmon (row 11)
$rec.Range("C11").Value = 4

# M.Alie-Cox (row 14)
$rec.Range("C14").Value = 32
$rec.Range("D14").Value = 19
$rec.Range("G14").Value = 10
$rec.Range("H14").Value = 6
